$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused FANN/NARX columns (B and C) entirely, for both rows,
# clearing their contents and formatting so the used range shrinks to A1:A2.
$ws.Range("B1:C2").Clear()

# Replace the remaining header/value with the new zip_code field.
$ws.Range("A1").Value = "zip_code"
$ws.Range("A2").Value = 94553

# Match the saved view's active selection.
$ws.Range("A2").Select()
